$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used to copy cell formatting (style) from a
# neighboring data row onto newly inserted rows.
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Remove the three rows that "fell off" the bottom of the tracked
#    window. Delete from the bottom-most row number up, so that row
#    numbers of not-yet-processed deletions remain valid.
# ---------------------------------------------------------------------
$ws.Rows.Item(40).Delete()
$ws.Rows.Item(37).Delete()
$ws.Rows.Item(19).Delete()

# ---------------------------------------------------------------------
# 2) Insert one new row (was "after" row 10) just above what is
#    currently row 8 (the old "+557388375022" row, untouched by the
#    deletes above since they all happened below it).
# ---------------------------------------------------------------------
$ws.Rows.Item(8).Insert()
$ws.Range("A9:C9").Copy()
$ws.Range("A8:C8").PasteSpecial($xlPasteFormats)
# Force text storage (so phone numbers / DDDs / dates are not
# reinterpreted as numbers or dates), then re-apply the plain data-row
# formatting so the final style index matches the other data rows.
$ws.Range("A8:C8").NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "+5511994806816"
$ws.Cells.Item(8, 2).Value = "11"
$ws.Cells.Item(8, 3).Value = "2024-10-23"
$ws.Range("A9:C9").Copy()
$ws.Range("A8:C8").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 3) Insert two brand-new rows at the very top of the data (rows 2-3),
#    pushing everything else down.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).Insert()
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial($xlPasteFormats)
$ws.Range("A2:C2").NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "+5516992621004"
$ws.Cells.Item(2, 2).Value = "16"
$ws.Cells.Item(2, 3).Value = "2024-10-24"
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial($xlPasteFormats)

$ws.Rows.Item(3).Insert()
$ws.Range("A4:C4").Copy()
$ws.Range("A3:C3").PasteSpecial($xlPasteFormats)
$ws.Range("A3:C3").NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "+5517996026345"
$ws.Cells.Item(3, 2).Value = "17"
$ws.Cells.Item(3, 3).Value = "2024-10-24"
$ws.Range("A4:C4").Copy()
$ws.Range("A3:C3").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false
